$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-05-26 Monday" "2025-05-27 Tuesday"
Replace-Text "412×2=" "679×3="
Replace-Text "206×8=" "338×9="
Replace-Text "324×9=" "309×2="
Replace-Text "171×4=" "536×6="
Replace-Text "915×6=" "199×6="
Replace-Text "644×9=" "580×9="
Replace-Text "643×7=" "725×9="
Replace-Text "886×5=" "430×6="
Replace-Text "752×7=" "460×7="
Replace-Text "318×3=" "767×9="
Replace-Text "500×5=" "350×2="
Replace-Text "917×2=" "724×6="
Replace-Text "681×6=" "894×3="
Replace-Text "738×2=" "755×4="
Replace-Text "548×4=" "659×7="
Replace-Text "163×3=" "226×8="
Replace-Text "363×9=" "424×4="
Replace-Text "880×6=" "294×6="
Replace-Text "412×4=" "765×8="
Replace-Text "961×2=" "758×7="
Replace-Text "291×8=" "919×5="
Replace-Text "473×7=" "167×4="
Replace-Text "873×2=" "435×3="
Replace-Text "914×9=" "717×7="
Replace-Text "217×9=" "232×2="
